$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-data refresh: it overwrites the
# cached price/profit columns (H, I, J, K, L, M, N) for specific rows across
# all eight crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# All of these columns hold static values (no formulas in the workbook), so
# each cell is written directly with its refreshed value.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4116313.8
$ws.Range("J112").Value = 4116313.8
$ws.Range("L112").Value = 12348941.4
$ws.Range("N112").Value = -12351157.4
$ws.Range("H116").Value = 4295.3335
$ws.Range("I116").Value = 1940.8
$ws.Range("K116").Value = 1940.8
$ws.Range("M116").Value = 1501.2
$ws.Range("H129").Value = 257613.39
$ws.Range("I129").Value = 797
$ws.Range("K129").Value = 2391
$ws.Range("M129").Value = 2609
$ws.Range("H132").Value = 2696.353
$ws.Range("I132").Value = 2839
$ws.Range("J132").Value = 1626.5
$ws.Range("K132").Value = 8517
$ws.Range("L132").Value = 4879.5
$ws.Range("M132").Value = -5987
$ws.Range("N132").Value = -9939.5
$ws.Range("H135").Value = 13161660
$ws.Range("I135").Value = 709.871
$ws.Range("J135").Value = 71445864
$ws.Range("K135").Value = 6388.839
$ws.Range("L135").Value = 643012776
$ws.Range("M135").Value = -3853.839
$ws.Range("N135").Value = -643017846
$ws.Range("H137").Value = 1566.3182
$ws.Range("I137").Value = 1213.25
$ws.Range("J137").Value = 1990
$ws.Range("K137").Value = 3639.75
$ws.Range("L137").Value = 5970
$ws.Range("M137").Value = -1089.75
$ws.Range("N137").Value = -11070
$ws.Range("H138").Value = 12823382
$ws.Range("I138").Value = 55557290
$ws.Range("J138").Value = 3209.75
$ws.Range("K138").Value = 166671870
$ws.Range("L138").Value = 9629.25
$ws.Range("M138").Value = -166666730
$ws.Range("N138").Value = -19909.25
$ws.Range("H141").Value = 1217.1305
$ws.Range("I141").Value = 859.8421
$ws.Range("J141").Value = 2914.25
$ws.Range("K141").Value = 2579.5263
$ws.Range("L141").Value = 8742.75
$ws.Range("M141").Value = 2600.4737
$ws.Range("N141").Value = -19102.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4315.884
$ws.Range("I32").Value = 4116.091
$ws.Range("K32").Value = 4116.091
$ws.Range("M32").Value = -3829.091
$ws.Range("H45").Value = 2646.6287
$ws.Range("I45").Value = 2400.28
$ws.Range("J45").Value = 3262.5
$ws.Range("K45").Value = 2400.28
$ws.Range("L45").Value = 3262.5
$ws.Range("M45").Value = -2023.28
$ws.Range("N45").Value = -4016.5
$ws.Range("H61").Value = 401003.47
$ws.Range("I61").Value = 487543.3
$ws.Range("J61").Value = 756.75
$ws.Range("K61").Value = 487543.3
$ws.Range("L61").Value = 756.75
$ws.Range("M61").Value = -487331.3
$ws.Range("N61").Value = -1180.75
$ws.Range("H132").Value = 11086.432
$ws.Range("I132").Value = 1248.3721
$ws.Range("J132").Value = 63966
$ws.Range("K132").Value = 3745.1163
$ws.Range("L132").Value = 191898
$ws.Range("M132").Value = -1215.1163
$ws.Range("N132").Value = -196958
$ws.Range("H136").Value = 401003.47
$ws.Range("I136").Value = 487543.3
$ws.Range("J136").Value = 756.75
$ws.Range("K136").Value = 1462629.9
$ws.Range("L136").Value = 2270.25
$ws.Range("M136").Value = -1460079.9
$ws.Range("N136").Value = -7370.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 11193
$ws.Range("J81").Value = 11193
$ws.Range("L81").Value = 11193
$ws.Range("N81").Value = -13315
$ws.Range("H84").Value = 11193
$ws.Range("J84").Value = 11193
$ws.Range("L84").Value = 33579
$ws.Range("N84").Value = -44187
$ws.Range("H105").Value = 4228.9
$ws.Range("I105").Value = 5177.8
$ws.Range("K105").Value = 5177.8
$ws.Range("M105").Value = -3430.8
$ws.Range("H134").Value = 4802.926
$ws.Range("I134").Value = 5835.737
$ws.Range("J134").Value = 2350
$ws.Range("K134").Value = 17507.211
$ws.Range("L134").Value = 7050
$ws.Range("M134").Value = -14972.211
$ws.Range("N134").Value = -12120
$ws.Range("H135").Value = 47695
$ws.Range("J135").Value = 47695
$ws.Range("L135").Value = 47695
$ws.Range("N135").Value = -57835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1974.6538
$ws.Range("I31").Value = 1194.1428
$ws.Range("J31").Value = 5252.8
$ws.Range("K31").Value = 1194.1428
$ws.Range("L31").Value = 5252.8
$ws.Range("M31").Value = -899.1428000000001
$ws.Range("N31").Value = -5842.8
$ws.Range("H34").Value = 1974.6538
$ws.Range("I34").Value = 1194.1428
$ws.Range("J34").Value = 5252.8
$ws.Range("K34").Value = 1194.1428
$ws.Range("L34").Value = 5252.8
$ws.Range("M34").Value = -992.1428000000001
$ws.Range("N34").Value = -5656.8
$ws.Range("H58").Value = 15587.471
$ws.Range("I58").Value = 927.931
$ws.Range("J58").Value = 100612.8
$ws.Range("K58").Value = 927.931
$ws.Range("L58").Value = 100612.8
$ws.Range("M58").Value = -724.931
$ws.Range("N58").Value = -101018.8
$ws.Range("H132").Value = 1829.711
$ws.Range("I132").Value = 1416.814
$ws.Range("J132").Value = 10707
$ws.Range("K132").Value = 4250.442
$ws.Range("L132").Value = 32121
$ws.Range("M132").Value = -1720.442
$ws.Range("N132").Value = -37181
$ws.Range("H134").Value = 802.9474
$ws.Range("I134").Value = 708.54346
$ws.Range("J134").Value = 1197.7273
$ws.Range("K134").Value = 2125.63038
$ws.Range("L134").Value = 3593.1819
$ws.Range("M134").Value = 409.3696199999999
$ws.Range("N134").Value = -8663.1819
$ws.Range("H136").Value = 15587.471
$ws.Range("I136").Value = 927.931
$ws.Range("J136").Value = 100612.8
$ws.Range("K136").Value = 2783.793
$ws.Range("L136").Value = 301838.4
$ws.Range("M136").Value = -233.7930000000001
$ws.Range("N136").Value = -306938.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 358.6
$ws.Range("J23").Value = 560.8333
$ws.Range("L23").Value = 1682.4999
$ws.Range("N23").Value = -2152.4999
$ws.Range("H131").Value = 703.76
$ws.Range("J131").Value = 717.1613
$ws.Range("L131").Value = 2151.4839
$ws.Range("N131").Value = -12231.4839

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 80.85714
$ws.Range("I2").Value = 81.71429
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 81.71429
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = 31.28570999999999
$ws.Range("N2").Value = -306
$ws.Range("H126").Value = 5271.4165
$ws.Range("I126").Value = 4281.25
$ws.Range("J126").Value = 7251.75
$ws.Range("K126").Value = 12843.75
$ws.Range("L126").Value = 21755.25
$ws.Range("M126").Value = -10373.75
$ws.Range("N126").Value = -26695.25
$ws.Range("H132").Value = 26710.092
$ws.Range("I132").Value = 4306.1
$ws.Range("K132").Value = 12918.3
$ws.Range("M132").Value = -10388.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47624372
$ws.Range("I7").Value = 83336504
$ws.Range("J7").Value = 8189.4443
$ws.Range("K7").Value = 83336504
$ws.Range("L7").Value = 8189.4443
$ws.Range("M7").Value = -83336392
$ws.Range("N7").Value = -8413.4443
$ws.Range("H46").Value = 774
$ws.Range("I46").Value = 840.8571
$ws.Range("J46").Value = 680.4
$ws.Range("K46").Value = 840.8571
$ws.Range("L46").Value = 680.4
$ws.Range("M46").Value = -652.8571
$ws.Range("N46").Value = -1056.4
$ws.Range("H126").Value = 47624372
$ws.Range("I126").Value = 83336504
$ws.Range("J126").Value = 8189.4443
$ws.Range("K126").Value = 250009512
$ws.Range("L126").Value = 24568.3329
$ws.Range("M126").Value = -250007042
$ws.Range("N126").Value = -29508.3329
$ws.Range("H132").Value = 345257.06
$ws.Range("I132").Value = 377578.66
$ws.Range("J132").Value = 493.33334
$ws.Range("K132").Value = 1132735.98
$ws.Range("L132").Value = 1480.00002
$ws.Range("M132").Value = -1130205.98
$ws.Range("N132").Value = -6540.000019999999
$ws.Range("H136").Value = 1061.5588
$ws.Range("I136").Value = 978.5517
$ws.Range("J136").Value = 1543
$ws.Range("K136").Value = 2935.6551
$ws.Range("L136").Value = 4629
$ws.Range("M136").Value = -385.6550999999999
$ws.Range("N136").Value = -9729

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 22185
$ws.Range("J104").Value = 22185
$ws.Range("L104").Value = 22185
$ws.Range("N104").Value = -29173
$ws.Range("H119").Value = 35000
$ws.Range("J119").Value = 35000
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676
$ws.Range("H126").Value = 1186.08
$ws.Range("I126").Value = 1131.3334
$ws.Range("K126").Value = 3394.0002
$ws.Range("M126").Value = -924.0002
$ws.Range("H132").Value = 425.58823
$ws.Range("I132").Value = 431.6129
$ws.Range("J132").Value = 363.33334
$ws.Range("K132").Value = 1294.8387
$ws.Range("L132").Value = 1090.00002
$ws.Range("M132").Value = 1235.1613
$ws.Range("N132").Value = -6150.000019999999
$ws.Range("I136").Value = 23461314
$ws.Range("J136").Value = 5336.6665
$ws.Range("K136").Value = 70383942
$ws.Range("L136").Value = 16009.9995
$ws.Range("M136").Value = -70383942
$ws.Range("N136").Value = -21109.9995
